$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.131.95"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.821.76"
$ws.Range("E3").Value = "  +1.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.46"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5106"
$ws.Range("E7").Value = "  -2.26%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3909"
$ws.Range("E8").Value = "  +3.04%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09715"
$ws.Range("E9").Value = "  +22.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.109"
$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.05"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.459"
$ws.Range("E12").Value = "  +3.44%  "

$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.005"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.85"
$ws.Range("E14").Value = "  +1.88%  "

$ws.Range("D15").Value = "1.827.81"
$ws.Range("E15").Value = "  +2.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.381"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001135"
$ws.Range("E17").Value = "  +4.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.04"
$ws.Range("E18").Value = "  +2.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06605"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.040"
$ws.Range("E22").Value = "  +1.36%  "

$ws.Range("D23").Value = "28.209.86"
$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.249"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.47"
$ws.Range("E26").Value = "  -1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.77"
$ws.Range("E27").Value = "  +1.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.446"
$ws.Range("E28").Value = "  +4.96%  "

$ws.Range("D29").Value = "2.038.51"
$ws.Range("E29").Value = "  +2.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.69"
$ws.Range("E30").Value = "  +4.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1091"
$ws.Range("E31").Value = "  +0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.065"
$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.638"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.638"
$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06917"
$ws.Range("E35").Value = "  -3.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.081"
$ws.Range("E36").Value = "  +5.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02342"
$ws.Range("E37").Value = "  +1.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2171"
$ws.Range("E38").Value = "  +1.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.59"
$ws.Range("E39").Value = "  -4.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.013"
$ws.Range("E40").Value = "  -0.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6254"
$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.155"
$ws.Range("E43").Value = "  -0.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.30"
$ws.Range("E44").Value = "  +0.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5972"
$ws.Range("E45").Value = "  +1.07%  "

$ws.Range("E46").Value = "  -1.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.289"
$ws.Range("E47").Value = "  -6.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.16"
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.993"
$ws.Range("E49").Value = "  +4.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.188"
$ws.Range("E50").Value = "  -1.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06792"
$ws.Range("E51").Value = "  +0.19%  "

